# Apply updates to SwaadSutra_Daily_2026-01-21.xlsx
$wb = $excel.ActiveWorkbook

# Order #27 (row 3) in "Daily Orders" moved from READY to COOKING status
$wsOrders = $wb.Worksheets.Item("Daily Orders")
$wsOrders.Range("H3").Value = "COOKING"

# Reflect the status change in the "Summary" sheet counts:
# Cooking count increases from 0 to 1, Ready count decreases from 1 to 0
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("C2").Value = 1
$wsSummary.Range("D2").Value = 0
